$d = $word.ActiveDocument

# Replace the instructional placeholder text with the "last modified by" notice.
$d.Content.Find.Execute(
    "Edit this document and submit a pull request for your instructor to review your changes.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This document has been modified by Gus Manley at 1:32 PM on December 2, 2015.",
    2
)
